$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.032.61'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').Value = '1.676.61'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3660'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '47.25'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3244'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.151'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07202'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9991'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.099'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '1.673.38'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.671'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001052'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06527'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9990'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '79.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.926'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.84%  '
$ws.Range('D24').Value = '25.019.65'
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.442'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.387'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').Value = '1.857.49'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.15'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.089'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.826'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08479'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.669'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.39'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.179'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06099'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.235'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.57%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02237'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.52%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2095'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.300'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.75%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5987'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.835'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5739'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.970'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07021'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.191'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.18%  '
